# Applies the scheduled market-data refresh captured in the commit diff.
# For each Leve row below we update the price/profit columns (H:N) to the
# latest scraped values. A few rows also gain or lose a trailing column
# (N = LeveProfitHQ) because the underlying scrape now has/lacks an HQ
# price quote for that item.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 6876.1875
$ws.Range("I33").Value = 245.8
$ws.Range("J33").Value = 9890
$ws.Range("K33").Value = 245.8
$ws.Range("L33").Value = 9890
$ws.Range("M33").Value = -16.80000000000001
$ws.Range("N33").Value = -10348

$ws.Range("H127").Value = 568.7692
$ws.Range("I127").Value = 354.8889
$ws.Range("K127").Value = 1064.6667
$ws.Range("M127").Value = 3895.3333

$ws.Range("H129").Value = 1024.8837
$ws.Range("I129").Value = 346.25
$ws.Range("J129").Value = 1180
$ws.Range("K129").Value = 1038.75
$ws.Range("L129").Value = 3540
$ws.Range("M129").Value = 3961.25
$ws.Range("N129").Value = -13540

$ws.Range("H131").Value = 1258.1818
$ws.Range("I131").Value = 704.44446
$ws.Range("J131").Value = 3750
$ws.Range("K131").Value = 2113.33338
$ws.Range("L131").Value = 11250
$ws.Range("M131").Value = 2926.66662
$ws.Range("N131").Value = -21330

$ws.Range("H132").Value = 3971.5144
$ws.Range("I132").Value = 4632.231
$ws.Range("J132").Value = 2062.7778
$ws.Range("K132").Value = 13896.693
$ws.Range("L132").Value = 6188.3334
$ws.Range("M132").Value = -11366.693
$ws.Range("N132").Value = -11248.3334

$ws.Range("H137").Value = 1794.25
$ws.Range("I137").Value = 1887.5264
$ws.Range("J137").Value = 1439.8
$ws.Range("K137").Value = 5662.5792
$ws.Range("L137").Value = 4319.4
$ws.Range("M137").Value = -3112.5792
$ws.Range("N137").Value = -9419.4

$ws.Range("H141").Value = 890.5161000000001
$ws.Range("I141").Value = 786.4138
$ws.Range("J141").Value = 2400
$ws.Range("K141").Value = 2359.2414
$ws.Range("L141").Value = 7200
$ws.Range("M141").Value = 2820.7586
$ws.Range("N141").Value = -17560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 29000.125
$ws.Range("J55").Value = 29000.125
$ws.Range("L55").Value = 29000.125
$ws.Range("N55").Value = -29630.125

$ws.Range("H74").Value = 773.6
$ws.Range("I74").Value = 677.2308
$ws.Range("J74").Value = 1400
$ws.Range("K74").Value = 677.2308
$ws.Range("L74").Value = 1400
$ws.Range("M74").Value = 196.7692
$ws.Range("N74").Value = -3148

$ws.Range("H77").Value = 773.6
$ws.Range("I77").Value = 677.2308
$ws.Range("J77").Value = 1400
$ws.Range("K77").Value = 3386.154
$ws.Range("L77").Value = 7000
$ws.Range("M77").Value = 981.8459999999995
$ws.Range("N77").Value = -15736

$ws.Range("H132").Value = 43523280
$ws.Range("I132").Value = 58824524
$ws.Range("J132").Value = 169752
$ws.Range("K132").Value = 176473572
$ws.Range("L132").Value = 509256
$ws.Range("M132").Value = -176471042
$ws.Range("N132").Value = -514316

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7355066.5
$ws.Range("I31").Value = 8622251
$ws.Range("K31").Value = 8622251
$ws.Range("M31").Value = -8621956

$ws.Range("H34").Value = 7355066.5
$ws.Range("I34").Value = 8622251
$ws.Range("K34").Value = 8622251
$ws.Range("M34").Value = -8622049

$ws.Range("H50").Value = 9857.429
$ws.Range("J50").Value = 9857.429
$ws.Range("L50").Value = 9857.429
$ws.Range("N50").Value = -11107.429

$ws.Range("H58").Value = 1144.1
$ws.Range("I58").Value = 1163.8182
$ws.Range("J58").Value = 1120
$ws.Range("K58").Value = 1163.8182
$ws.Range("L58").Value = 1120
$ws.Range("M58").Value = -960.8181999999999
$ws.Range("N58").Value = -1526

$ws.Range("H59").Value = 12645.125
$ws.Range("J59").Value = 12645.125
$ws.Range("L59").Value = 12645.125
$ws.Range("N59").Value = -14935.125

$ws.Range("H60").Value = 7667
$ws.Range("J60").Value = 10000.5
$ws.Range("L60").Value = 10000.5
$ws.Range("N60").Value = -11022.5

$ws.Range("H68").Value = 17600.4
$ws.Range("J68").Value = 17600.4
$ws.Range("L68").Value = 17600.4
$ws.Range("N68").Value = -19098.4

$ws.Range("H71").Value = 17600.4
$ws.Range("J71").Value = 17600.4
$ws.Range("L71").Value = 52801.2
$ws.Range("N71").Value = -60289.2

$ws.Range("H74").Value = 19240.666
$ws.Range("J74").Value = 19240.666
$ws.Range("L74").Value = 19240.666
$ws.Range("N74").Value = -20988.666

$ws.Range("H77").Value = 19240.666
$ws.Range("J77").Value = 19240.666
$ws.Range("L77").Value = 57721.99800000001
$ws.Range("N77").Value = -66457.99800000001

$ws.Range("H80").Value = 20966.666
$ws.Range("J80").Value = 21450
$ws.Range("L80").Value = 21450
$ws.Range("N80").Value = -23696

$ws.Range("H83").Value = 20966.666
$ws.Range("J83").Value = 21450
$ws.Range("L83").Value = 64350
$ws.Range("N83").Value = -75582

$ws.Range("H134").Value = 5558570.5
$ws.Range("I134").Value = 6668375.5
$ws.Range("J134").Value = 9543.333000000001
$ws.Range("K134").Value = 20005126.5
$ws.Range("L134").Value = 28629.999
$ws.Range("M134").Value = -20002591.5
$ws.Range("N134").Value = -33699.999

$ws.Range("H136").Value = 1144.1
$ws.Range("I136").Value = 1163.8182
$ws.Range("J136").Value = 1120
$ws.Range("K136").Value = 3491.4546
$ws.Range("L136").Value = 3360
$ws.Range("M136").Value = -941.4546
$ws.Range("N136").Value = -8460

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1583.3334

$ws.Range("H55").Value = 3110
$ws.Range("I55").Value = 626.6667
$ws.Range("J55").Value = 4600
$ws.Range("K55").Value = 1880.0001
$ws.Range("L55").Value = 13800
$ws.Range("M55").Value = -1703.0001
$ws.Range("N55").Value = -14154

$ws.Range("H80").Value = 2185
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2185
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 6555
$ws.Range("N80").Value = -8427
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 2185
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2185
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 19665
$ws.Range("N83").Value = -29025
$ws.Range("M83").ClearContents()

$ws.Range("H135").Value = 1583.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 12000
$ws.Range("J62").Value = 12000
$ws.Range("L62").Value = 12000
$ws.Range("N62").Value = -13248

$ws.Range("H65").Value = 12000
$ws.Range("J65").Value = 12000
$ws.Range("L65").Value = 36000
$ws.Range("N65").Value = -42240

$ws.Range("H132").Value = 91245.414
$ws.Range("I132").Value = 120374
$ws.Range("J132").Value = 3859.6667
$ws.Range("K132").Value = 361122
$ws.Range("L132").Value = 11579.0001
$ws.Range("M132").Value = -358592
$ws.Range("N132").Value = -16639.0001

$ws.Range("H136").Value = 15645.538
$ws.Range("I136").Value = 17938
$ws.Range("J136").Value = 11977.6
$ws.Range("K136").Value = 53814
$ws.Range("L136").Value = 35932.8
$ws.Range("M136").Value = -51264
$ws.Range("N136").Value = -41032.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H136").Value = 84763.164
$ws.Range("I136").Value = 126057.25
$ws.Range("J136").Value = 2175
$ws.Range("K136").Value = 378171.75
$ws.Range("L136").Value = 6525
$ws.Range("M136").Value = -375621.75
$ws.Range("N136").Value = -11625
